$p = $ppt.ActivePresentation

# Slide 17 ("Data Structures Introduction") - body placeholder (shape 2) resized/repositioned
$s17 = $p.Slides.Item(17)
$sh17 = $s17.Shapes.Item(2)
$sh17.Left = 12.5
$sh17.Top = 191.96418002834645
$sh17.Width = 695.0
$sh17.Height = 156.0716535433071

# Slide 21 ("Java Collection Framework") - body placeholder (shape 2) resized/repositioned
$s21 = $p.Slides.Item(21)
$sh21 = $s21.Shapes.Item(2)
$sh21.Left = 12.5
$sh21.Top = 182.4396133992126
$sh21.Width = 695.0
$sh21.Height = 175.1207874015748

# Slide 25 ("Maps") - body placeholder (shape 2) resized/repositioned
$s25 = $p.Slides.Item(25)
$sh25 = $s25.Shapes.Item(2)
$sh25.Left = 12.43527559055118
$sh25.Top = 187.7212601425197
$sh25.Width = 695.0
$sh25.Height = 164.55748751496063
